$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Apply new custom number format (0;[Red]0) to the "Box No" header and
#        the existing count column (C1:C11). This creates numFmtId 164 plus
#        two new cellXfs entries (indices 7 and 8) mirroring the old ones.
$ws.Range("C1").NumberFormat = "0;[Red]0"
$ws.Range("C2:C11").NumberFormat = "0;[Red]0"

# --- 2. Fix a handful of existing values that were corrected in this edit.
$ws.Cells.Item(3, 1).Value = 24061691
$ws.Cells.Item(6, 1).Value = 24061612
$ws.Cells.Item(11, 1).Value = 2406145

$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(4, 3).Value = 3
$ws.Cells.Item(5, 3).Value = 4
$ws.Cells.Item(6, 3).Value = 5
$ws.Cells.Item(7, 3).Value = 6
$ws.Cells.Item(8, 3).Value = 7
$ws.Cells.Item(9, 3).Value = 8
$ws.Cells.Item(10, 3).Value = 9
$ws.Cells.Item(11, 3).Value = 10

# --- 3. Append new rows 12-51 with body/cover box numbers.
$newRows = @(
    @(12, 2406146,  "C0624126", 11),
    @(13, 24061463, "C0624055", 12),
    @(14, 2406147,  "C0624076", 13),
    @(15, 2406148,  "C0624001", 14),
    @(16, 2406149,  "C0624030", 15),
    @(17, 2406150,  "C0624009", 16),
    @(18, 2406151,  "C0624077", 17),
    @(19, 2406152,  "C 0524082", 18),
    @(20, 2406153,  "C 0624129", 19),
    @(21, 2406154,  "C 0424079", 20),
    @(22, 2406155,  "C0624126", 21),
    @(23, 2406156,  "C0624055", 22),
    @(24, 2406157,  "C0624076", 23),
    @(25, 2406158,  "C0624001", 24),
    @(26, 2406159,  "C0624030", 25),
    @(27, 2406160,  "C0624009", 26),
    @(28, 2406161,  "C0624077", 27),
    @(29, 2406162,  "C 0524082", 28),
    @(30, 2406163,  "C 0624129", 29),
    @(31, 2406164,  "C 0424079", 30),
    @(32, 2406165,  "C0624126", 31),
    @(33, 2406166,  "C0624055", 32),
    @(34, 2406167,  "C0624076", 33),
    @(35, 2406168,  "C0624001", 34),
    @(36, 2406169,  "C0624030", 35),
    @(37, 2406170,  "C0624009", 36),
    @(38, 2406171,  "C0624077", 37),
    @(39, 2406172,  "C 0524082", 38),
    @(40, 2406173,  "C 0624129", 39),
    @(41, 2406174,  "C 0424079", 40),
    @(42, 2406175,  "C0624126", 41),
    @(43, 2406176,  "C0624055", 42),
    @(44, 2406177,  "C0624076", 43),
    @(45, 2406178,  "C0624001", 44),
    @(46, 2406179,  "C0624030", 45),
    @(47, 2406180,  "C0624009", 46),
    @(48, 2406181,  "C0624077", 47),
    @(49, 2406182,  "C 0524082", 48),
    @(50, 2406183,  "C 0624129", 49)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
}

# Row 51 only has a value in column C.
$ws.Cells.Item(51, 3).Value = 50

# --- 4. Give the new B-column cells (B12:B50) the same font/border/wrap
#        formatting as the rest of the B column (style index 1), matching
#        B2:B11.
$ws.Range("B11").Copy()
$ws.Range("B12:B50").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 5. Update the active selection to match the saved view state.
$ws.Range("G6").Select()
